# Generate Report for Handoff
# - Set Priority ("ht") for the six "Ready for handoff" rows (7-12) on both
#   the zh-cn and de-de localization-status sheets.
# - Refresh the "Latest Handoff Datetime" timestamps for those same rows on
#   each sheet (the Overview sheet's "Latest HO Xliff Generate Date" shares
#   the de-de timestamp text and updates along with it).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Priority column (E) for rows 7-12: mark as handed off ("ht").
$zhcn.Range("E7:E12").Value = "ht"
$dede.Range("E7:E12").Value = "ht"

# Latest Handoff Datetime column (H) for rows 7-12.
$zhcn.Range("H7:H12").Value = "2016-08-18 02:18:10"
$dede.Range("H7:H12").Value = "2016-08-18 02:18:15"

# Overview sheet mirrors the de-de handoff timestamp in its
# "Latest HO Xliff Generate Date" column (G) for the same rows.
$overview.Range("G7:G12").Value = "2016-08-18 02:18:15"
